$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Handback report generation: append a new row (row 4) to each of the three
# sheets (Overview, zh-cn, de-de) describing the hand-back of
# 7d100aed-79da-4ebd-a115-59c39d97faed.md, mirroring the existing rows.
# ---------------------------------------------------------------------------

$commitSrc   = "9f1a2b3c4d5e6f708192a3b4c5d6e7f8091a2b3c"
$commitZhCn  = "b2c3d4e5f60718293a4b5c6d7e8f90a1b2c3d4e5"
$commitDeDe  = "c3d4e5f60718293a4b5c6d7e8f90a1b2c3d4e5f6"

$mdFile   = "7d100aed-79da-4ebd-a115-59c39d97faed.md"
$mdPath   = "e2e\7d100aed-79da-4ebd-a115-59c39d97faed.md"
$xlfZhCn  = "7d100aed-79da-4ebd-a115-59c39d97faed.27d3d8f644be3618c6059a79b181ec26831a1a73.zh-cn.xlf"
$xlfDeDe  = "7d100aed-79da-4ebd-a115-59c39d97faed.27d3d8f644be3618c6059a79b181ec26831a1a73.de-de.xlf"

$statusInSync = "Handed back: in sync with en-US"

$dtHandoff   = "2016-08-31 04:45:54"
$dtZhHo      = "2016-08-31 04:45:49"
$dtZhHb      = "2016-08-31 04:46:14"
$dtDeHb      = "2016-08-31 04:46:21"

$dateStyleFormat = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------------
# Sheet "Overview" -> new row 4
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

$wsOv.Range("A4").Value = $mdFile
$wsOv.Range("C4").Value = ".md"
$wsOv.Range("E4").Value = $statusInSync
$wsOv.Range("F4").Value = $statusInSync
$wsOv.Range("G4").Value = $dtHandoff
$wsOv.Range("G4").NumberFormat = $dateStyleFormat

$wsOv.Hyperlinks.Add($wsOv.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSrc/e2e/$mdFile", "", "", $mdPath)

# ---------------------------------------------------------------------------
# Sheet "zh-cn" -> new row 4
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B4").Value = ".md"
$wsZh.Range("C4").Value = $statusInSync
$wsZh.Range("D4").Value = "e2e"
$wsZh.Range("E4").Value = "ht"
$wsZh.Range("F4").Value = "True"
$wsZh.Range("G4").Value = $xlfZhCn
$wsZh.Range("H4").Value = $dtZhHo
$wsZh.Range("H4").NumberFormat = $dateStyleFormat
$wsZh.Range("J4").Value = $xlfZhCn
$wsZh.Range("K4").Value = $dtZhHb
$wsZh.Range("K4").NumberFormat = $dateStyleFormat
$wsZh.Range("L4").Value = ""
$wsZh.Range("M4").Value = "True"
$wsZh.Range("N4").Value = ""
$wsZh.Range("O4").Value = "False"
$wsZh.Range("P4").Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/$commitZhCn/e2e/$mdFile", "", "", $mdFile)
$wsZh.Hyperlinks.Add($wsZh.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/$commitZhCn/e2e/$mdFile", "", "", $mdFile)

# ---------------------------------------------------------------------------
# Sheet "de-de" -> new row 4
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B4").Value = ".md"
$wsDe.Range("C4").Value = $statusInSync
$wsDe.Range("D4").Value = "e2e"
$wsDe.Range("E4").Value = "ht"
$wsDe.Range("F4").Value = "True"
$wsDe.Range("G4").Value = $xlfDeDe
$wsDe.Range("H4").Value = $dtHandoff
$wsDe.Range("H4").NumberFormat = $dateStyleFormat
$wsDe.Range("J4").Value = $xlfDeDe
$wsDe.Range("K4").Value = $dtDeHb
$wsDe.Range("K4").NumberFormat = $dateStyleFormat
$wsDe.Range("L4").Value = ""
$wsDe.Range("M4").Value = "True"
$wsDe.Range("N4").Value = ""
$wsDe.Range("O4").Value = "False"
$wsDe.Range("P4").Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/$commitDeDe/e2e/$mdFile", "", "", $mdFile)
$wsDe.Hyperlinks.Add($wsDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/$commitDeDe/e2e/$mdFile", "", "", $mdFile)

# ---------------------------------------------------------------------------
# Keep the tables / autofilters in sync with the newly added row.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("Overview").ListObjects.Item("Overview").Resize($wsOv.Range("A1:G4"))
$wb.Worksheets.Item("zh-cn").ListObjects.Item("zh_cn").Resize($wsZh.Range("A1:P4"))
$wb.Worksheets.Item("de-de").ListObjects.Item("de_de").Resize($wsDe.Range("A1:P4"))
